$d = $word.ActiveDocument

# The document contains six "<id>...</id>" tags, each currently split across
# three separate runs (opening "<id>" tag run, the inner identifier text run,
# and the closing "</id>" tag run). This collapses each one back into a
# single run (keeping the Courier New / 7f6000 formatting used by the
# tag-delimiter runs) by finding the whole "<id>...</id>" span and replacing
# it in-place with itself - Word's Find/Replace merges the matched span into
# a single run using the formatting of its first character run.

$ids = @("p039v_3", "p040r_1", "p040r_2", "p040r_3", "p040r_4", "p040r_5")

foreach ($id in $ids) {
    $needle = "<id>" + $id + "</id>"
    $d.Content.Find.Execute($needle, $true, $false, $false, $false, $false, $true, 1, $false, $needle, 2) | Out-Null
}
